$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new CNPJ/test data ("lojas que estão dando erro")
$ws.Range("A2").Value = "93.421.457/0001-50"
$ws.Range("B2").Value = "Teste"
$ws.Range("C2").Value = "Teste"
$ws.Range("B2:C2").Font.Underline = $true

# Remove the old sample rows (row 3 previously held TesteCNPJ01 data)
$ws.Range("A3:C3").ClearContents()

# Build out the array area (rows 3 through 22) to hold the stores with errors
$ws.Range("A3:C22").Font.Underline = $true

$ws.Range("C2").Select()
